$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Mumbaya / magnetic sensor broken ---
$ws.Range("B2").Value = "'0848"
$ws.Range("C2").Value = "Mumbaya"
$ws.Range("D2").Value = "Setor 1 magnético do cliente está quebrado, tem de trocar."
$ws.Range("E2").ClearContents()
$ws.Range("G2").Value = "Pendente"

# --- Row 3: Igreja Batista Betel / disparo zona 4, battery fixed ---
$ws.Range("B3").Value = "'0682"
$ws.Range("C3").Value = "Igreja Batista Betel"
$ws.Range("D3").Value = "Disparo frequente na zona 4."
$ws.Range("E3").Value = "Foi trocado a bateria do sensor e foi revisado o sistema."
$ws.Range("F3").Value = "Bateria do setor foi trocada."
$ws.Range("G3").Value = "Concluido"

# --- Row 4: Rc Silva stays, description updated ---
$ws.Range("D4").Value = "Zona aberta, cliente pedindo reparo."

# --- Row 5: RotoPlast stays, description updated ---
$ws.Range("D5").Value = "Central acusando falha de rede elétrica e câmera fora."

# --- Row 6: Escola São Geraldo / disparo zona 13 ---
$ws.Range("B6").Value = "'0081"
$ws.Range("C6").Value = "Escola São Geraldo"
$ws.Range("D6").Value = "Disparo na zona 13 após o arme."
$ws.Range("E6").ClearContents()
$ws.Range("G6").Value = "Pendente"
$ws.Rows(6).AutoFit()

# --- Row 7: Depósito Ideal / sem comunicação de alarmes ---
$ws.Range("B7").Value = "'0210"
$ws.Range("C7").Value = "Depósito Ideal"
$ws.Range("D7").Value = "Sem comunicação de alarmes, linha telefônica."
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = "Pendente"

# --- Row 8: clear everything except the I column note ---
$ws.Range("A8").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Rows(8).AutoFit()

# --- Rows 9-16: clear entirely ---
$ws.Range("A9:I16").ClearContents()
$ws.Rows(9).AutoFit()
$ws.Rows(15).AutoFit()

# --- Update selection to H7 ---
$ws.Range("H7").Select()
